# Applies the "Updated cryptos list" data refresh to the cryptocurrency
# table on the active worksheet (prices in column D, 1h volume % in
# column E). Row 37/38 additionally swap the coin identity (ImmutableX
# and VeChain traded places in the ranking).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into column D (Price) while forcing it to be
# stored as text, matching the original inline-string cell type -
# otherwise Excel would happily reinterpret strings like "0.503" or
# "214.14" as numbers (losing the original formatting / precision).
function Set-Text($rangeAddr, $val) {
    $cell = $ws.Range($rangeAddr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-Text "D2" '26.450.53'
$ws.Range("E2").Value = '  +0.87%  '
Set-Text "D3" '1.617.18'
$ws.Range("E3").Value = '  +1.80%  '
$ws.Range("E4").Value = '  -0.21%  '
Set-Text "D5" '214.14'
$ws.Range("E5").Value = '  +0.91%  '
Set-Text "D6" '0.503'
$ws.Range("E6").Value = '  +0.43%  '
$ws.Range("E7").Value = '  -0.16%  '
Set-Text "D8" '0.247'
$ws.Range("E8").Value = '  +0.33%  '
$ws.Range("E9").Value = '  +0.43%  '
Set-Text "D10" '19.24'
$ws.Range("E10").Value = '  -0.57%  '
Set-Text "D11" '0.0856'
$ws.Range("E11").Value = '  +0.87%  '
Set-Text "D12" '1.844.57'
$ws.Range("E12").Value = '  +1.84%  '
Set-Text "D13" '1.606.56'
$ws.Range("E13").Value = '  +0.29%  '
$ws.Range("E14").Value = '  +0.02%  '
$ws.Range("E15").Value = '  -1.46%  '
Set-Text "D16" '64.73'
$ws.Range("E16").Value = '  +0.57%  '
Set-Text "D17" '26.458.82'
$ws.Range("E17").Value = '  +0.86%  '
Set-Text "D18" '229.76'
$ws.Range("E18").Value = '  +7.70%  '
$ws.Range("E19").Value = '  +0.05%  '
Set-Text "D20" '7.59'
$ws.Range("E20").Value = '  +2.59%  '
Set-Text "D21" '0.999'
$ws.Range("E21").Value = '  -0.19%  '
$ws.Range("E22").Value = '  +1.97%  '
$ws.Range("E23").Value = '  +1.34%  '
$ws.Range("E24").Value = '  -0.46%  '
Set-Text "D25" '145.32'
$ws.Range("E25").Value = '  +0.93%  '
$ws.Range("E26").Value = '  -0.12%  '
Set-Text "D27" '7.03'
$ws.Range("E27").Value = '  -0.35%  '
$ws.Range("E28").Value = '  +2.38%  '
Set-Text "D29" '15.56'
$ws.Range("E29").Value = '  +2.44%  '
$ws.Range("E30").Value = '  +0.14%  '
$ws.Range("E31").Value = '  +0.58%  '
$ws.Range("E32").Value = '  +0.85%  '
Set-Text "D33" '1.452.75'
$ws.Range("E33").Value = '  +9.03%  '
$ws.Range("E34").Value = '  +2.19%  '
$ws.Range("E35").Value = '  -0.87%  '
$ws.Range("E36").Value = '  +0.47%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-Text "D37" '0.557'
$ws.Range("E37").Value = '  -5.39%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-Text "D38" '0.0167'
$ws.Range("E38").Value = '  +0.35%  '
Set-Text "D39" '0.835'
$ws.Range("E39").Value = '  +2.34%  '
Set-Text "D40" '5.87'
$ws.Range("E40").Value = '  +2.45%  '
$ws.Range("E41").Value = '  -0.11%  '
Set-Text "D42" '2.19'
$ws.Range("E42").Value = '  +2.36%  '
Set-Text "D43" '1.756.01'
$ws.Range("E43").Value = '  +1.89%  '
$ws.Range("E44").Value = '  -0.34%  '
$ws.Range("E45").Value = '  +0.33%  '
Set-Text "D46" '0.914'
$ws.Range("E46").Value = '  -9.70%  '
Set-Text "D47" '88.04'
$ws.Range("E47").Value = '  +3.02%  '
$ws.Range("E48").Value = '  +0.56%  '
Set-Text "D49" '0.0502'
$ws.Range("E49").Value = '  +0.11%  '
$ws.Range("E50").Value = '  -1.11%  '
Set-Text "D51" '7.47'
$ws.Range("E51").Value = '  +1.53%  '

Write-Host "cryptos list updated"
